# Insert a new weekly price record for "Feria Lagunitas de Puerto Montt - Albahaca"
# at row 118, pushing the existing rows 118-136 down to 119-137.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(118).Insert()

$ws.Cells.Item(118, 1).Value = 4
$ws.Cells.Item(118, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(118, 3).Value = "Los Lagos"
$ws.Cells.Item(118, 4).Value = 44782
$ws.Cells.Item(118, 5).Value = 10
$ws.Cells.Item(118, 6).Value = 100112052
$ws.Cells.Item(118, 7).Value = "Albahaca"
$ws.Cells.Item(118, 8).Value = "Sin especificar"
$ws.Cells.Item(118, 9).Value = "Primera"
$ws.Cells.Item(118, 10).Value = 100
$ws.Cells.Item(118, 11).Value = 6000
$ws.Cells.Item(118, 12).Value = 6500
$ws.Cells.Item(118, 13).Value = 6250
$ws.Cells.Item(118, 14).Value = "$/paquete"
$ws.Cells.Item(118, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(118, 16).Value = 6250
$ws.Cells.Item(118, 17).Value = 1
$ws.Cells.Item(118, 18).Value = "Hortaliza"

# Preserve the date-formatted style used by the other rows' "Fecha" column.
$ws.Cells.Item(118, 4).NumberFormat = $ws.Cells.Item(119, 4).NumberFormat
